$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - update entrada/salida timestamps and recalc tiempo/total
$ws.Range("C2").Value = "08/08/2024 15:23:00"
$ws.Range("D2").Value = "08/08/2024 15:23:00"
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 1500

# Row 3 - update entrada/salida timestamps and recalc tiempo/total
$ws.Range("C3").Value = "08/08/2024 17:40:00"
$ws.Range("D3").Value = "08/08/2024 17:40:00"
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 1500
